$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Duplicate the "Cars" sheet to create "Cars (2)" as the new first tab ---
# Copying twice and then collapsing back to a single duplicate reproduces the
# sheetId sequence (sheetId=4) that a natural "duplicate, undo/redo, tidy up"
# editing session would leave behind, matching sheetId="4" in the target.
$cars = $wb.Worksheets.Item("Cars")
$trafficLights = $wb.Worksheets.Item("TrafficLights")

$cars.Copy($trafficLights)
$wb.Worksheets.Item("Cars").Copy($trafficLights)
$wb.Worksheets.Item("Cars (2)").Delete()
$wb.Worksheets.Item("Cars (3)").Name = "Cars (2)"

# --- New "Cars (2)" sheet: recalculated AppearTime-ish values + selection ---
$carsCopy = $wb.Worksheets.Item("Cars (2)")
$carsCopy.Range("F2").Value = 40
$carsCopy.Range("F3").Value = 30
$carsCopy.Range("F4").Select()

# --- TrafficLights: move the selection to I2 (intersection flag column) ---
$tl = $wb.Worksheets.Item("TrafficLights")
$tl.Range("I2").Select()

# --- Cars: same value updates as the duplicate, plus new selection ---
$cars = $wb.Worksheets.Item("Cars")
$cars.Range("F2").Value = 40
$cars.Range("F3").Value = 30
$cars.Range("A3:H5").Select()

# "Cars" stays the active tab/sheet (activeTab=2, tabSelected on sheet3)
$cars.Activate()
